$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.450759333333334
$ws.Range("H2").Value = 28.352278
$ws.Range("I2").Value = 0.6100581096599255
$ws.Range("J2").Value = 0.6100581096599255
$ws.Range("M2").Value = 75.02398799999999
$ws.Range("N2").Value = 225.071964
$ws.Range("O2").Value = 0.9946207163417996
$ws.Range("P2").Value = 0.9946207163417997
$ws.Range("Q2").Value = 709.0336548148879
$ws.Range("R2").Value = 6381.302893333992
$ws.Range("S2").Value = 0.6067764340400792
$ws.Range("T2").Value = 0.6067764340400793
$ws.Range("G3").Value = 9.450759333333334
$ws.Range("H3").Value = 28.352278
$ws.Range("I3").Value = 0.6100581096599255
$ws.Range("J3").Value = 0.6100581096599255
$ws.Range("O3").Value = 0.001281651759329359
$ws.Range("P3").Value = 0.001281651759329359
$ws.Range("Q3").Value = 0.913649008296889
$ws.Range("R3").Value = 8.222841074672
$ws.Range("S3").Value = 0.0007818820495387864
$ws.Range("T3").Value = 0.0007818820495387864
$ws.Range("G4").Value = 9.450759333333334
$ws.Range("H4").Value = 28.352278
$ws.Range("I4").Value = 0.6100581096599255
$ws.Range("J4").Value = 0.6100581096599255
$ws.Range("M4").Value = 0.291865
$ws.Range("N4").Value = 0.875595
$ws.Range("O4").Value = 0.00386936209489556
$ws.Range("P4").Value = 0.00386936209489556
$ws.Range("Q4").Value = 2.758345872823333
$ws.Range("R4").Value = 24.82511285541
$ws.Range("S4").Value = 0.002360535725201755
$ws.Range("T4").Value = 0.002360535725201755
$ws.Range("G5").Value = 9.450759333333334
$ws.Range("H5").Value = 28.352278
$ws.Range("I5").Value = 0.6100581096599255
$ws.Range("J5").Value = 0.6100581096599255
$ws.Range("M5").Value = 0.01721833333333333
$ws.Range("N5").Value = 0.051655
$ws.Range("O5").Value = 0.0002282698039753883
$ws.Range("P5").Value = 0.0002282698039753884
$ws.Range("Q5").Value = 0.1627263244544444
$ws.Range("R5").Value = 1.46453692009
$ws.Range("S5").Value = 0.0001392578451056671
$ws.Range("T5").Value = 0.0001392578451056672
$ws.Range("H6").Value = 5.309089
$ws.Range("I6").Value = 0.1142360694740755
$ws.Range("J6").Value = 0.1142360694740756
$ws.Range("M6").Value = 75.02398799999999
$ws.Range("N6").Value = 225.071964
$ws.Range("O6").Value = 0.9946207163417996
$ws.Range("P6").Value = 0.9946207163417997
$ws.Range("Q6").Value = 132.769676475644
$ws.Range("R6").Value = 1194.927088280796
$ws.Range("S6").Value = 0.1136215612523766
$ws.Range("T6").Value = 0.1136215612523766
$ws.Range("H7").Value = 5.309089
$ws.Range("I7").Value = 0.1142360694740755
$ws.Range("J7").Value = 0.1142360694740756
$ws.Range("O7").Value = 0.001281651759329359
$ws.Range("P7").Value = 0.001281651759329359
$ws.Range("S7").Value = 0.0001464108594203198
$ws.Range("T7").Value = 0.0001464108594203198
$ws.Range("H8").Value = 5.309089
$ws.Range("I8").Value = 0.1142360694740755
$ws.Range("J8").Value = 0.1142360694740756
$ws.Range("M8").Value = 0.291865
$ws.Range("N8").Value = 0.875595
$ws.Range("O8").Value = 0.00386936209489556
$ws.Range("P8").Value = 0.00386936209489556
$ws.Range("Q8").Value = 0.5165124203283333
$ws.Range("R8").Value = 4.648611782955
$ws.Range("S8").Value = 0.0004420207170928437
$ws.Range("T8").Value = 0.0004420207170928438
$ws.Range("H9").Value = 5.309089
$ws.Range("I9").Value = 0.1142360694740755
$ws.Range("J9").Value = 0.1142360694740756
$ws.Range("M9").Value = 0.01721833333333333
$ws.Range("N9").Value = 0.051655
$ws.Range("O9").Value = 0.0002282698039753883
$ws.Range("P9").Value = 0.0002282698039753884
$ws.Range("Q9").Value = 0.03047122136611111
$ws.Range("R9").Value = 0.274240992295
$ws.Range("S9").Value = 0.00002607664518576607
$ws.Range("T9").Value = 0.00002607664518576607
$ws.Range("G10").Value = 3.398621333333333
$ws.Range("H10").Value = 10.195864
$ws.Range("I10").Value = 0.2193851766757396
$ws.Range("J10").Value = 0.2193851766757396
$ws.Range("M10").Value = 75.02398799999999
$ws.Range("N10").Value = 225.071964
$ws.Range("O10").Value = 0.9946207163417996
$ws.Range("P10").Value = 0.9946207163417997
$ws.Range("Q10").Value = 254.9781261285439
$ws.Range("R10").Value = 2294.803135156896
$ws.Range("S10").Value = 0.2182050415799964
$ws.Range("T10").Value = 0.2182050415799965
$ws.Range("G11").Value = 3.398621333333333
$ws.Range("H11").Value = 10.195864
$ws.Range("I11").Value = 0.2193851766757396
$ws.Range("J11").Value = 0.2193851766757396
$ws.Range("O11").Value = 0.001281651759329359
$ws.Range("P11").Value = 0.001281651759329359
$ws.Range("Q11").Value = 0.3285605845262222
$ws.Range("R11").Value = 2.957045260736
$ws.Range("S11").Value = 0.0002811753976572439
$ws.Range("T11").Value = 0.0002811753976572439
$ws.Range("G12").Value = 3.398621333333333
$ws.Range("H12").Value = 10.195864
$ws.Range("I12").Value = 0.2193851766757396
$ws.Range("J12").Value = 0.2193851766757396
$ws.Range("M12").Value = 0.291865
$ws.Range("N12").Value = 0.875595
$ws.Range("O12").Value = 0.00386936209489556
$ws.Range("P12").Value = 0.00386936209489556
$ws.Range("Q12").Value = 0.9919386154533333
$ws.Range("R12").Value = 8.927447539080001
$ws.Range("S12").Value = 0.0008488806868110724
$ws.Range("T12").Value = 0.0008488806868110726
$ws.Range("G13").Value = 3.398621333333333
$ws.Range("H13").Value = 10.195864
$ws.Range("I13").Value = 0.2193851766757396
$ws.Range("J13").Value = 0.2193851766757396
$ws.Range("M13").Value = 0.01721833333333333
$ws.Range("N13").Value = 0.051655
$ws.Range("O13").Value = 0.0002282698039753883
$ws.Range("P13").Value = 0.0002282698039753884
$ws.Range("Q13").Value = 0.0585185949911111
$ws.Range("R13").Value = 0.52666735492
$ws.Range("S13").Value = 0.00005007901127487702
$ws.Range("T13").Value = 0.00005007901127487703
$ws.Range("G14").Value = 0.8724953333333333
$ws.Range("H14").Value = 2.617486
$ws.Range("I14").Value = 0.0563206441902594
$ws.Range("J14").Value = 0.05632064419025941
$ws.Range("M14").Value = 75.02398799999999
$ws.Range("N14").Value = 225.071964
$ws.Range("O14").Value = 0.9946207163417996
$ws.Range("P14").Value = 0.9946207163417997
$ws.Range("Q14").Value = 65.45807941805599
$ws.Range("R14").Value = 589.122714762504
$ws.Range("S14").Value = 0.05601767946934742
$ws.Range("T14").Value = 0.05601767946934744
$ws.Range("G15").Value = 0.8724953333333333
$ws.Range("H15").Value = 2.617486
$ws.Range("I15").Value = 0.0563206441902594
$ws.Range("J15").Value = 0.05632064419025941
$ws.Range("O15").Value = 0.001281651759329359
$ws.Range("P15").Value = 0.001281651759329359
$ws.Range("Q15").Value = 0.08434819551822223
$ws.Range("R15").Value = 0.759133759664
$ws.Range("S15").Value = 0.0000721834527130088
$ws.Range("T15").Value = 0.00007218345271300881
$ws.Range("G16").Value = 0.8724953333333333
$ws.Range("H16").Value = 2.617486
$ws.Range("I16").Value = 0.0563206441902594
$ws.Range("J16").Value = 0.05632064419025941
$ws.Range("M16").Value = 0.291865
$ws.Range("N16").Value = 0.875595
$ws.Range("O16").Value = 0.00386936209489556
$ws.Range("P16").Value = 0.00386936209489556
$ws.Range("Q16").Value = 0.2546508504633333
$ws.Range("R16").Value = 2.29185765417
$ws.Range("S16").Value = 0.0002179249657898896
$ws.Range("T16").Value = 0.0002179249657898896
$ws.Range("G17").Value = 0.8724953333333333
$ws.Range("H17").Value = 2.617486
$ws.Range("I17").Value = 0.0563206441902594
$ws.Range("J17").Value = 0.05632064419025941
$ws.Range("M17").Value = 0.01721833333333333
$ws.Range("N17").Value = 0.051655
$ws.Range("O17").Value = 0.0002282698039753883
$ws.Range("P17").Value = 0.0002282698039753884
$ws.Range("Q17").Value = 0.015022915481111105
$ws.Range("R17").Value = 0.13520623933
$ws.Range("S17").Value = 0.000012856302409078106
$ws.Range("T17").Value = 0.000012856302409078112
